$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.732.69'
$ws.Range('E2').Value = '  -0.81%  '
$ws.Range('D3').Value = '2.353.61'
$ws.Range('E3').Value = '  -4.06%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '542.42'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.99%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.44'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.80%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.523'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -11.25%  '
$ws.Range('D9').Value = '2.352.00'
$ws.Range('E9').Value = '  -4.06%  '
$ws.Range('E10').Value = '  -2.48%  '
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('E12').Value = '  -3.62%  '
$ws.Range('E13').Value = '  -3.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '24.81'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.65%  '
$ws.Range('D15').Value = '2.775.66'
$ws.Range('E15').Value = '  -4.13%  '
$ws.Range('D16').Value = '60.510.26'
$ws.Range('E16').Value = '  -1.07%  '
$ws.Range('E17').Value = '  -3.00%  '
$ws.Range('D18').Value = '2.352.20'
$ws.Range('E18').Value = '  -4.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.63'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.95%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '314.44'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.05%  '
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.06'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.86%  '
$ws.Range('E22').Value = '  -7.93%  '
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.88'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.21'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.83%  '
$ws.Range('E26').Value = '  +7.81%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('D28').Value = '2.468.02'
$ws.Range('E28').Value = '  -4.57%  '
$ws.Range('D29').Value = '0.0₃0893'
$ws.Range('E29').Value = '  -8.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.95'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.78%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.38'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.52%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '500.36'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -8.65%  '
$ws.Range('E33').Value = '  -1.65%  '
$ws.Range('E34').Value = '  -5.51%  '
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.55'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.53%  '
$ws.Range('E38').Value = '  -1.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.41'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.22'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -10.85%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.80'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.84%  '
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '138.93'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.07%  '
$ws.Range('B43').Value = 'USDe'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.12'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.72%  '
$ws.Range('E45').Value = '  -10.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '138.53'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.01%  '
$ws.Range('E47').Value = '  -1.71%  '
$ws.Range('E48').Value = '  -4.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.46'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -8.89%  '
$ws.Range('E50').Value = '  -3.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0894'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.13%  '
